$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.352.16"
$ws.Range("E2").Value = "  -1.61%  "

$ws.Range("D3").Value = "1.707.92"
$ws.Range("E3").Value = "  -1.49%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "223.90"
$ws.Range("E5").Value = "  -1.20%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5302"
$ws.Range("E6").Value = "  -2.05%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.005"
$ws.Range("E7").Value = "  +0.16%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2668"
$ws.Range("E8").Value = "  -2.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06603"
$ws.Range("E9").Value = "  -1.33%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.78"
$ws.Range("E10").Value = "  -5.32%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07613"
$ws.Range("E11").Value = "  -1.65%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.533"
$ws.Range("E12").Value = "  -3.03%  "

$ws.Range("D13").Value = "1.726.48"
$ws.Range("E13").Value = "  -0.57%  "

$ws.Range("D14").Value = "1.944.82"
$ws.Range("E14").Value = "  -1.42%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5759"
$ws.Range("E15").Value = "  -2.92%  "

$ws.Range("D16").Value = "0.0₅8155"
$ws.Range("E16").Value = "  -2.75%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.57"
$ws.Range("E17").Value = "  -1.70%  "

$ws.Range("D18").Value = "27.353.47"
$ws.Range("E18").Value = "  -1.62%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "216.38"
$ws.Range("E19").Value = "  -4.56%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.004"
$ws.Range("E20").Value = "  +0.13%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.645"
$ws.Range("E21").Value = "  -3.26%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.42"
$ws.Range("E22").Value = "  -3.85%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.939"
$ws.Range("E23").Value = "  -4.22%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.005"
$ws.Range("E24").Value = "  +0.12%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "142.05"
$ws.Range("E25").Value = "  -4.12%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.716"
$ws.Range("E26").Value = "  -0.58%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1205"
$ws.Range("E27").Value = "  -3.14%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.232"
$ws.Range("E28").Value = "  -3.40%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "16.13"
$ws.Range("E29").Value = "  -5.20%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05385"
$ws.Range("E30").Value = "  -4.24%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.288"
$ws.Range("E31").Value = "  -1.71%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.475"
$ws.Range("E32").Value = "  -4.86%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.416"
$ws.Range("E33").Value = "  -2.41%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.642"
$ws.Range("E34").Value = "  -2.30%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.878"
$ws.Range("E35").Value = "  +0.97%  "

$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.412"
$ws.Range("E36").Value = "  -1.13%  "

$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9443"
$ws.Range("E37").Value = "  -2.85%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5832"
$ws.Range("E38").Value = "  -2.43%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01634"
$ws.Range("E39").Value = "  -2.34%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.785"
$ws.Range("E40").Value = "  -2.49%  "

$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.005"
$ws.Range("E41").Value = "  +0.16%  "

$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "1.043.36"
$ws.Range("E42").Value = "  -1.26%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8409"
$ws.Range("E43").Value = "  -2.89%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.98"
$ws.Range("E44").Value = "  -0.77%  "

$ws.Range("D45").Value = "1.851.98"
$ws.Range("E45").Value = "  -1.36%  "

$ws.Range("E46").Value = "  +2.76%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "57.75"
$ws.Range("E47").Value = "  -3.46%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4527"
$ws.Range("E48").Value = "  +2.13%  "

$ws.Range("E49").Value = "  +0.19%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.086"
$ws.Range("E50").Value = "  -2.35%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05248"
$ws.Range("E51").Value = "  -1.00%  "
